# Actualización automática 2025-08-25 16:45:08
#
# Applies updated figures across the three worksheets of the workbook:
#   - "VENTAS POR GRUPO"
#   - "VENTA MENSUAL"
#   - "CUMPLIMIENTO MENSUAL"

$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasPorGrupo.Range("D4").Value = 950.4
$wsVentasPorGrupo.Range("M27").Value = 1147.45
$wsVentasPorGrupo.Range("M34").Value = 9635.41
$wsVentasPorGrupo.Range("D55").Value = "11 de 53"

# --- Sheet: VENTA MENSUAL ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F4").Value = 1323.71
$wsVentaMensual.Range("F27").Value = 3886.33
$wsVentaMensual.Range("F34").Value = 16519.07
$wsVentaMensual.Range("F55").Value = 65109.4

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumplimientoMensual = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumplimientoMensual.Range("D3").Value = 17782.84
$wsCumplimientoMensual.Range("E3").Value = 9674.167600000001
$wsCumplimientoMensual.Range("F3").Value = 0.6476612549723008

$wsCumplimientoMensual.Range("D16").Value = 33399.97
$wsCumplimientoMensual.Range("E16").Value = 22659.73
$wsCumplimientoMensual.Range("F16").Value = 0.5957928779497572

$wsCumplimientoMensual.Range("D19").Value = 71856.60000000001
$wsCumplimientoMensual.Range("E19").Value = 45583.09064517915
$wsCumplimientoMensual.Range("F19").Value = 0.611859581758441
